$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Clear first so the old shared string becomes fully unreferenced (and is
# dropped) before the new text is written.
$ws1.Range("E2:F3").ClearContents()
$ws2.Range("C2:C3").ClearContents()
$ws3.Range("C2:C3").ClearContents()

$ws1.Range("E2:F3").Value = "In Translation"
$ws2.Range("C2:C3").Value = "In Translation"
$ws3.Range("C2:C3").Value = "In Translation"

# --- Column width changes (17.2159881591797 -> 13.4101845877511 character units) ---
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
